$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04628467895650446
$ws.Range("C2").Value = 1.329827146587384
$ws.Range("D2").Value = 9.907280297060181
$ws.Range("E2").Value = 3.147583247042114
$ws.Range("F2").Value = 3.177948460597182
$ws.Range("G2").Value = 52

$ws.Range("B3").Value = 0.03072812890343799
$ws.Range("C3").Value = 1.487177437779083
$ws.Range("D3").Value = 9.489403539255518
$ws.Range("E3").Value = 3.080487548953172
$ws.Range("F3").Value = 3.110985134611961
$ws.Range("G3").Value = 51

$ws.Range("B4").Value = 0.05670300065403502
$ws.Range("C4").Value = 1.350516681017773
$ws.Range("D4").Value = 6.724916598423166
$ws.Range("E4").Value = 2.593244415480956
$ws.Range("F4").Value = 2.61894615200115
$ws.Range("G4").Value = 50

$ws.Range("B5").Value = 0.06807619608839131
$ws.Range("C5").Value = 1.406782146406506
$ws.Range("D5").Value = 8.32307609620424
$ws.Range("E5").Value = 2.884974193334187
$ws.Range("F5").Value = 2.914059470405175
$ws.Range("G5").Value = 49

$ws.Range("B6").Value = 0.0465829568975624
$ws.Range("C6").Value = 1.552000945598137
$ws.Range("D6").Value = 9.764321243549908
$ws.Range("E6").Value = 3.124791392005218
$ws.Range("F6").Value = 3.15750797590164
$ws.Range("G6").Value = 48

$ws.Range("B7").Value = 0.007348863100462911
$ws.Range("C7").Value = 1.800235720797896
$ws.Range("D7").Value = 10.52634790614585
$ws.Range("E7").Value = 3.244433372123066
$ws.Range("F7").Value = 3.290447561450098
$ws.Range("G7").Value = 36

$ws.Range("B8").Value = 0.04734274986515632
$ws.Range("C8").Value = 1.848121295394991
$ws.Range("D8").Value = 10.96964521407898
$ws.Range("E8").Value = 3.312045472827778
$ws.Range("F8").Value = 3.360055740670333
$ws.Range("G8").Value = 35

$ws.Range("B9").Value = 0.08958096046880724
$ws.Range("C9").Value = 2.597962177260555
$ws.Range("D9").Value = 19.61373083157788
$ws.Range("E9").Value = 4.428739192092698
$ws.Range("F9").Value = 4.55620268551751
$ws.Range("G9").Value = 18

$ws.Range("B10").Value = -0.6075801358359439
$ws.Range("C10").Value = 3.106648192679968
$ws.Range("D10").Value = 23.3436798468976
$ws.Range("E10").Value = 4.831529762600827
$ws.Range("F10").Value = 5.027124311967897
$ws.Range("G10").Value = 11
